$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C684").Value = "RÏSɒTÖ"
$ws.Range("C686").Value = "RɒS©"
$ws.Range("C687").Value = "RɒSIÑ"
$ws.Range("C688").Value = "Rɒθ©"
$ws.Range("C689").Value = "RɒθS"
$ws.Range("C690").Value = "RÚDɒLF/RÚDɒF"
$ws.Range("C691").Value = "RʌNɒF"
$ws.Range("C692").Value = "SɒNKÚLɒT"
$ws.Range("C693").Value = "SɛRəTɒV/SɛRəTɒF"
$ws.Range("C694").Value = "S©ɒÑ"
$ws.Range("C695").Value = "S©ɒÑZ"
$ws.Range("C698").Value = "SKɒF"
$ws.Range("C699").Value = "SKɒFђ"
$ws.Range("C700").Value = "SKɒFS"
$ws.Range("C703").Value = "SNɒT"
$ws.Range("C704").Value = "SNɒTFÉS"
$ws.Range("C705").Value = "SNɒTFÉSђ"
$ws.Range("C706").Value = "SNɒTÏ"
$ws.Range("C707").Value = "SɒDəM"
$ws.Range("C708").Value = "SɒFIT"
$ws.Range("C710").Value = "SɒFTBɔL/SɒFBɔL"
$ws.Range("C711").Value = "SɒFTBɔLZ/SɒFBɔLZ"
$ws.Range("C713").Value = "SɒFəNћ"
$ws.Range("C714").Value = "SɒFəN©"
$ws.Range("C715").Value = "SɒFəNIÑ/SɒFNIÑ"
$ws.Range("C716").Value = "SɒFəNZ"
$ws.Range("C718").Value = "SɒFTəST"
$ws.Range("C719").Value = "SɒFTLÝ/SɒFLÝ"
$ws.Range("C720").Value = "SɒFTNəS/SɒFNəS"
$ws.Range("C721").Value = "SɒFTSPÖKəN/SɒFSPÖKəN"
$ws.Range("C722").Value = "SɒFTWɛR/SɒFWɛR"
$ws.Range("C723").Value = "SɒFTWɛRZ/SɒFWɛRZ"
$ws.Range("C724").Value = "SɒFTWUD"
$ws.Range("C725").Value = "SɒLəMLÝ"
$ws.Range("C726").Value = "SɒLSTIS"
$ws.Range("C727").Value = "SɒLSTISIS"
$ws.Range("C728").Value = "SɒLVəNSÝ"
$ws.Range("C729").Value = "SɒLVəNTS"
$ws.Range("C730").Value = "SɒÑ"
$ws.Range("C731").Value = "SɒÑZ"
$ws.Range("C732").Value = "SɒÑB®D"
$ws.Range("C733").Value = "SɒÑB®DZ"
$ws.Range("C734").Value = "SɒÑ©"
$ws.Range("C735").Value = "SɒÑRÍT©"
$ws.Range("C736").Value = "SɒÑRÍT©Z"
$ws.Range("C737").Value = "SɒÑRÍTIÑ"
$ws.Range("C738").Value = "SɒFMɒR"
$ws.Range("C739").Value = "SɒFMɒRZ"
$ws.Range("C743").Value = "SɒT"
$ws.Range("C744").Value = "SPÆZMɒDIK"
$ws.Range("C745").Value = "SPÆZMɒDIKLÝ"
$ws.Range("C747").Value = "STÆNDɒF"
$ws.Range("C748").Value = "STÆNDɒFS"
$ws.Range("C749").Value = "STRɛPTəKɒKəS"
$ws.Range("C750").Value = "STRɒÑ"
$ws.Range("C751").Value = "STRɒÑG©"
$ws.Range("C752").Value = "STRɒÑGəST"
$ws.Range("C753").Value = "STRɒÑXÖLD"
$ws.Range("C754").Value = "STRɒÑXÖLDZ"
$ws.Range("C755").Value = "STRɒÑLÝ"
$ws.Range("C760").Value = "SWɒMPђ"
$ws.Range("C763").Value = "SWɒN©"
$ws.Range("C766").Value = "SINəGɒG"
$ws.Range("C767").Value = "SINəGɒGZ"
$ws.Range("C768").Value = "TÉBəLKLɒθ"
$ws.Range("C769").Value = "TÉBəLKLɒθS"
$ws.Range("C770").Value = "TÉKɒF"
$ws.Range("C771").Value = "TÉKɒFS"
$ws.Range("C772").Value = "TÏLÝəLɒJIKəL"
$ws.Range("C773").Value = "TÚTɒNIK"
$ws.Range("C774").Value = "ÐɛRɒN"
$ws.Range("C775").Value = "θɒÑ"
$ws.Range("C776").Value = "θÁZəNDYIRZLɒÑ"
$ws.Range("C777").Value = "θRɒÑ"
$ws.Range("C778").Value = "θRɒÑћ"
$ws.Range("C779").Value = "θRɒÑZ"
$ws.Range("C780").Value = "TɒF"
$ws.Range("C781").Value = "TɒL©ÉTS"
$ws.Range("C784").Value = "TɒÑ"
$ws.Range("C787").Value = "TɒPDÁN"
$ws.Range("C788").Value = "TɒPSÏKRIT"
$ws.Range("C789").Value = "TəPɒLəJÝ"
$ws.Range("C796").Value = "TɒS"
$ws.Range("C797").Value = "TɒSђ"
$ws.Range("C798").Value = "TɒSIZ"
$ws.Range("C799").Value = "TɒSIÑ"
$ws.Range("C800").Value = "TRÉDɒF"
$ws.Range("C801").Value = "TRÉDɒFS"
$ws.Range("C802").Value = "TRÆVəLɒG"
$ws.Range("C803").Value = "TRÏTɒP"
$ws.Range("C804").Value = "TRÏTɒPS"
$ws.Range("C807").Value = "ʌNKɒNTRəDIKTʔћ"
$ws.Range("C808").Value = "ʌND©DɒG"
$ws.Range("C809").Value = "ʌND©DɒGZ"
$ws.Range("C810").Value = "ʌND©GɒN"
$ws.Range("C811").Value = "əNGɒDLÝ"
$ws.Range("C813").Value = "əNINVɒLVћ"
$ws.Range("C818").Value = "əNWɒNTʔћ"
$ws.Range("C822").Value = "VɒLəNTIRIÑ"
$ws.Range("C823").Value = "VɒLəNTIRIZəM"
$ws.Range("C825").Value = "VUTɒN"
$ws.Range("C829").Value = "WɒLəT"
$ws.Range("C830").Value = "WɒLəTS"
$ws.Range("C831").Value = "WɒLəPIÑ"
$ws.Range("C833").Value = "WɒNTʔћ"
$ws.Range("C837").Value = "WɔRLɒK"
$ws.Range("C838").Value = "WɔRLɒKS"
$ws.Range("C852").Value = "WɒʃKLɒθ"
$ws.Range("C856").Value = "WɒÇDɒG"
$ws.Range("C857").Value = "WɒÇDɒGZ"
$ws.Range("C859").Value = "WɔT©BɒTəL"
$ws.Range("C860").Value = "WɔT©BɒTəLZ"
$ws.Range("C863").Value = "WÖBIGɒN"
$ws.Range("C864").Value = "WULFDɒG"
$ws.Range("C865").Value = "RÍTɒF"
$ws.Range("C866").Value = "RÍTɒFS"
$ws.Range("C867").Value = "RɒÑDÚ©"
$ws.Range("C868").Value = "RɒÑDÚ©Z"
$ws.Range("C869").Value = "RɒÑDÚIÑ"
$ws.Range("C870").Value = "RɒÑDÚIÑZ"
$ws.Range("C871").Value = "RɒÑћ"
$ws.Range("C872").Value = "RɒÑFəL"
$ws.Range("C873").Value = "RɒÑFəLÝ"
$ws.Range("C874").Value = "RɒÑXɛDʔћ"
$ws.Range("C875").Value = "RɒÑLÝ"
$ws.Range("C876").Value = "RɒÑZ"
$ws.Range("C877").Value = "YIRLɒÑ"

$ws.Application.ActiveWindow.ScrollRow = 703
$ws.Range("C712").Select()
